# Applies the "fixed quantification design tables" edit to Supplementary Table S1.
# Strategy: the table gains a 7th "Notes" column (populated for the header and the
# three method rows), the footnote row is merged into one (gridSpan 6 -> 7) with new
# wording, the second footnote row is removed, the eimeriaSpecies counts are corrected,
# and the stray superscript "2" marker on "Amplicon sequencing" is dropped.
$d = $word.ActiveDocument
$newBodyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"><w:body>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
        <w:pStyle w:val="caption"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">Table </w:t>
      </w:r>
      <w:r>
        <w:fldChar w:fldCharType="begin" w:dirty="true"/>
      </w:r>
      <w:r>
        <w:instrText xml:space="preserve" w:dirty="true"> SEQ Table \* ARABIC </w:instrText>
      </w:r>
      <w:r>
        <w:fldChar w:fldCharType="separate" w:dirty="true"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="default">1</w:t>
      </w:r>
      <w:r>
        <w:fldChar w:fldCharType="end" w:dirty="true"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="24"/>
          <w:color w:val="333333"/>
        </w:rPr>
        <w:t xml:space="default">Eimeria detection methods in field mice</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
        <w:pStyle w:val="caption"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
          <w:color w:val="333333"/>
        </w:rPr>
        <w:t xml:space="default">Summary of hierarchical approach for infection status and species assignment</w:t>
      </w:r>
    </w:p><w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:tblPr><w:tblCellMar><w:top w:w="0" w:type="dxa"></w:top><w:bottom w:w="0" w:type="dxa"></w:bottom><w:start w:w="60" w:type="dxa"></w:start><w:end w:w="60" w:type="dxa"></w:end></w:tblCellMar><w:tblW w:type="auto" w:w="0"></w:tblW><w:tblLook w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"></w:tblLook><w:jc w:val="center"></w:jc></w:tblPr><w:tr><w:trPr><w:cantSplit></w:cantSplit><w:tblHeader></w:tblHeader></w:trPr><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start></w:tcBorders></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Detection Method</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="superscript"/>
          <w:i/>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">1</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:bottom></w:tcBorders></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Variable Name</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:bottom></w:tcBorders></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Purpose</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:bottom></w:tcBorders></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Priority</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:bottom></w:tcBorders></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Sample Size (n/total)</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:bottom></w:tcBorders></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Success Rate</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"></w:bottom><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Notes</w:t>
      </w:r>
    </w:p></w:tc></w:tr><w:tr><w:trPr><w:cantSplit></w:cantSplit></w:trPr><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Caecal qPCR + melting curve</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">MC.Eimeria</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Infection detection (presence/absence)</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Primary</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">185/336</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">55.1%</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Direct infection status from melting curve</w:t>
      </w:r>
    </w:p></w:tc></w:tr>
    <w:tr><w:trPr><w:cantSplit></w:cantSplit></w:trPr><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Caecal qPCR + melting curve</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">eimeriaSpecies</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Species identification</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Primary</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">49/336</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">14.6%</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="F3E5F5"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Species ID from melting curve patterns</w:t>
      </w:r>
    </w:p></w:tc></w:tr>
    <w:tr><w:trPr><w:cantSplit></w:cantSplit></w:trPr><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFF3E0"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Amplicon sequencing</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFF3E0"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">amplicon_species</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFF3E0"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Backup species identification</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFF3E0"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Backup</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFF3E0"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">134/336</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFF3E0"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="end"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">39.9%</w:t>
      </w:r>
    </w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:space="0" w:color="D3D3D3"></w:top><w:bottom w:val="single" w:space="0" w:color="D3D3D3"></w:bottom><w:start w:val="single" w:space="0" w:color="D3D3D3"></w:start><w:end w:val="single" w:space="0" w:color="D3D3D3"></w:end></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFF3E0"></w:shd></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
        <w:jc w:val="start"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Used when qPCR species ID unavailable</w:t>
      </w:r>
    </w:p></w:tc></w:tr><w:tr><w:trPr><w:cantSplit></w:cantSplit></w:trPr><w:tc><w:tcPr><w:gridSpan w:val="7"/></w:tcPr><w:p>
      <w:pPr>
        <w:spacing w:before="0" w:after="60"/>
        <w:keepNext/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="superscript"/>
          <w:i/>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="20"/>
        </w:rPr>
        <w:t xml:space="default">Final species assignment: 49 from qPCR + 120 from amplicon = 169 total mice</w:t>
      </w:r>
    </w:p></w:tc></w:tr></w:tbl>
    <w:p><w:pPr><w:pStyle w:val="FirstParagraph" /></w:pPr></w:p>
    <w:sectPr />
  </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($newBodyXml)
Write-Host "Table S1 updated: added Notes column, fixed footnote, corrected counts."
